$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.342589
$ws.Range("H2").Value = 16.027767
$ws.Range("I2").Value = 0.4438787133356475
$ws.Range("J2").Value = 0.4621170122195802
$ws.Range("M2").Value = 0.9817236666666668
$ws.Range("N2").Value = 2.945171
$ws.Range("O2").Value = 0.04688329326954743
$ws.Range("P2").Value = 0.04832841473263862
$ws.Range("Q2").Value = 5.244946062573001
$ws.Range("R2").Value = 47.20451456315701
$ws.Range("S2").Value = 0.02081049589342453
$ws.Range("T2").Value = 0.0223333826215557

$ws.Range("G3").Value = 5.342589
$ws.Range("H3").Value = 16.027767
$ws.Range("I3").Value = 0.4438787133356475
$ws.Range("J3").Value = 0.4621170122195802
$ws.Range("M3").Value = 5.154927333333333
$ws.Range("N3").Value = 15.464782
$ws.Range("O3").Value = 0.2461792235003055
$ws.Range("P3").Value = 0.2537674037418691
$ws.Range("Q3").Value = 27.540658066866
$ws.Range("R3").Value = 247.865922601794
$ws.Range("S3").Value = 0.1092737169772844
$ws.Range("T3").Value = 0.1172702344159125

$ws.Range("G4").Value = 5.342589
$ws.Range("H4").Value = 16.027767
$ws.Range("I4").Value = 0.4438787133356475
$ws.Range("J4").Value = 0.4621170122195802
$ws.Range("M4").Value = 1.8784265
$ws.Range("N4").Value = 3.756853
$ws.Range("O4").Value = 0.08970632314876403
$ws.Range("P4").Value = 0.06164760887349412
$ws.Range("Q4").Value = 10.0356607562085
$ws.Range("R4").Value = 60.21396453725099
$ws.Range("S4").Value = 0.03981872729734519
$ws.Range("T4").Value = 0.02848840882310038

$ws.Range("G5").Value = 5.342589
$ws.Range("H5").Value = 16.027767
$ws.Range("I5").Value = 0.4438787133356475
$ws.Range("J5").Value = 0.4621170122195802
$ws.Range("M5").Value = 12.924656
$ws.Range("N5").Value = 38.773968
$ws.Range("O5").Value = 0.617231160081383
$ws.Range("P5").Value = 0.6362565726519981
$ws.Range("Q5").Value = 69.051124974384
$ws.Range("R5").Value = 621.460124769456
$ws.Range("S5").Value = 0.2739757731675934
$ws.Range("T5").Value = 0.2940249863590116

$ws.Range("I6").Value = 0.3624119935622804
$ws.Range("J6").Value = 0.3773029492651114
$ws.Range("M6").Value = 0.9817236666666668
$ws.Range("N6").Value = 2.945171
$ws.Range("O6").Value = 0.04688329326954743
$ws.Range("P6").Value = 0.04832841473263862
$ws.Range("Q6").Value = 4.282321502600111
$ws.Range("R6").Value = 38.540893523401
$ws.Range("S6").Value = 0.01699106777858173
$ws.Range("T6").Value = 0.01823445341193201

$ws.Range("I7").Value = 0.3624119935622804
$ws.Range("J7").Value = 0.3773029492651114
$ws.Range("M7").Value = 5.154927333333333
$ws.Range("N7").Value = 15.464782
$ws.Range("O7").Value = 0.2461792235003055
$ws.Range("P7").Value = 0.2537674037418691
$ws.Range("Q7").Value = 22.48601812649355
$ws.Range("R7").Value = 202.374163138442
$ws.Range("S7").Value = 0.08921830316235989
$ws.Range("T7").Value = 0.09574718985915749

$ws.Range("I8").Value = 0.3624119935622804
$ws.Range("J8").Value = 0.3773029492651114
$ws.Range("M8").Value = 1.8784265
$ws.Range("N8").Value = 3.756853
$ws.Range("O8").Value = 0.08970632314876403
$ws.Range("P8").Value = 0.06164760887349412
$ws.Range("Q8").Value = 8.193778417623832
$ws.Range("R8").Value = 49.16267050574299
$ws.Range("S8").Value = 0.03251064740748572
$ws.Range("T8").Value = 0.02325982464311139

$ws.Range("I9").Value = 0.3624119935622804
$ws.Range("J9").Value = 0.3773029492651114
$ws.Range("M9").Value = 12.924656
$ws.Range("N9").Value = 38.773968
$ws.Range("O9").Value = 0.617231160081383
$ws.Range("P9").Value = 0.6362565726519981
$ws.Range("Q9").Value = 56.37791384864533
$ws.Range("R9").Value = 507.401224637808
$ws.Range("S9").Value = 0.223691975213853
$ws.Range("T9").Value = 0.2400614813509105

$ws.Range("G10").Value = 0.8673346666666665
$ws.Range("H10").Value = 2.602004
$ws.Range("I10").Value = 0.07206082965981524
$ws.Range("J10").Value = 0.07502169917140648
$ws.Range("M10").Value = 0.9817236666666668
$ws.Range("N10").Value = 2.945171
$ws.Range("O10").Value = 0.04688329326954743
$ws.Range("P10").Value = 0.04832841473263862
$ws.Range("Q10").Value = 0.851482969187111
$ws.Range("R10").Value = 7.663346722684
$ws.Range("S10").Value = 0.003378449010188019
$ws.Range("T10").Value = 0.003625679791502984

$ws.Range("G11").Value = 0.8673346666666665
$ws.Range("H11").Value = 2.602004
$ws.Range("I11").Value = 0.07206082965981524
$ws.Range("J11").Value = 0.07502169917140648
$ws.Range("M11").Value = 5.154927333333333
$ws.Range("N11").Value = 15.464782
$ws.Range("O11").Value = 0.2461792235003055
$ws.Range("P11").Value = 0.2537674037418691
$ws.Range("Q11").Value = 4.471047180347554
$ws.Range("R11").Value = 40.239424623128
$ws.Range("S11").Value = 0.0177398790904411
$ws.Range("T11").Value = 0.01903806182303136

$ws.Range("G12").Value = 0.8673346666666665
$ws.Range("H12").Value = 2.602004
$ws.Range("I12").Value = 0.07206082965981524
$ws.Range("J12").Value = 0.07502169917140648
$ws.Range("M12").Value = 1.8784265
$ws.Range("N12").Value = 3.756853
$ws.Range("O12").Value = 0.08970632314876403
$ws.Range("P12").Value = 0.06164760887349412
$ws.Range("Q12").Value = 1.629224422235333
$ws.Range("R12").Value = 9.775346533411998
$ws.Range("S12").Value = 0.006464312071831425
$ws.Range("T12").Value = 0.004624908367543805

$ws.Range("G13").Value = 0.8673346666666665
$ws.Range("H13").Value = 2.602004
$ws.Range("I13").Value = 0.07206082965981524
$ws.Range("J13").Value = 0.07502169917140648
$ws.Range("M13").Value = 12.924656
$ws.Range("N13").Value = 38.773968
$ws.Range("O13").Value = 0.617231160081383
$ws.Range("P13").Value = 0.6362565726519981
$ws.Range("Q13").Value = 11.21000220354133
$ws.Range("R13").Value = 100.890019831872
$ws.Range("S13").Value = 0.04447818948735469
$ws.Range("T13").Value = 0.04773304918932834

$ws.Range("G14").Value = 1.425086
$ws.Range("H14").Value = 2.850172
$ws.Range("I14").Value = 0.1184005245532914
$ws.Range("J14").Value = 0.08217694760298831
$ws.Range("M14").Value = 0.9817236666666668
$ws.Range("N14").Value = 2.945171
$ws.Range("O14").Value = 0.04688329326954743
$ws.Range("P14").Value = 0.04832841473263862
$ws.Range("Q14").Value = 1.399040653235333
$ws.Range("R14").Value = 8.394243919412
$ws.Range("S14").Value = 0.00555100651590021
$ws.Range("T14").Value = 0.003971481605219532

$ws.Range("G15").Value = 1.425086
$ws.Range("H15").Value = 2.850172
$ws.Range("I15").Value = 0.1184005245532914
$ws.Range("J15").Value = 0.08217694760298831
$ws.Range("M15").Value = 5.154927333333333
$ws.Range("N15").Value = 15.464782
$ws.Range("O15").Value = 0.2461792235003055
$ws.Range("P15").Value = 0.2537674037418691
$ws.Range("Q15").Value = 7.346214773750666
$ws.Range("R15").Value = 44.07728864250399
$ws.Range("S15").Value = 0.02914774919655812
$ws.Range("T15").Value = 0.02085383064064196

$ws.Range("G16").Value = 1.425086
$ws.Range("H16").Value = 2.850172
$ws.Range("I16").Value = 0.1184005245532914
$ws.Range("J16").Value = 0.08217694760298831
$ws.Range("M16").Value = 1.8784265
$ws.Range("N16").Value = 3.756853
$ws.Range("O16").Value = 0.08970632314876403
$ws.Range("P16").Value = 0.06164760887349412
$ws.Range("Q16").Value = 2.676919307178999
$ws.Range("R16").Value = 10.707677228716
$ws.Range("S16").Value = 0.01062127571656073
$ws.Range("T16").Value = 0.005066012324246643

$ws.Range("G17").Value = 1.425086
$ws.Range("H17").Value = 2.850172
$ws.Range("I17").Value = 0.1184005245532914
$ws.Range("J17").Value = 0.08217694760298831
$ws.Range("M17").Value = 12.924656
$ws.Range("N17").Value = 38.773968
$ws.Range("O17").Value = 0.617231160081383
$ws.Range("P17").Value = 0.6362565726519981
$ws.Range("Q17").Value = 18.41874632041599
$ws.Range("R17").Value = 110.512477922496
$ws.Range("S17").Value = 0.07308049312427231
$ws.Range("T17").Value = 0.05228562303288017

$ws.Range("G18").Value = 0.03909266666666666
$ws.Range("H18").Value = 0.117278
$ws.Range("I18").Value = 0.00324793888896551
$ws.Range("J18").Value = 0.003381391740913623
$ws.Range("M18").Value = 0.9817236666666668
$ws.Range("N18").Value = 2.945171
$ws.Range("O18").Value = 0.04688329326954743
$ws.Range("P18").Value = 0.04832841473263862
$ws.Range("Q18").Value = 0.03837819605977778
$ws.Range("R18").Value = 0.345403764538
$ws.Range("S18").Value = 0.000152274071452938
$ws.Range("T18").Value = 0.0001634173024283925

$ws.Range("G19").Value = 0.03909266666666666
$ws.Range("H19").Value = 0.117278
$ws.Range("I19").Value = 0.00324793888896551
$ws.Range("J19").Value = 0.003381391740913623
$ws.Range("M19").Value = 5.154927333333333
$ws.Range("N19").Value = 15.464782
$ws.Range("O19").Value = 0.2461792235003055
$ws.Range("P19").Value = 0.2537674037418691
$ws.Range("Q19").Value = 0.2015198559328889
$ws.Range("R19").Value = 1.813678703396
$ws.Range("S19").Value = 0.000799575073661974
$ws.Range("T19").Value = 0.0008580870031258491

$ws.Range("G20").Value = 0.03909266666666666
$ws.Range("H20").Value = 0.117278
$ws.Range("I20").Value = 0.00324793888896551
$ws.Range("J20").Value = 0.003381391740913623
$ws.Range("M20").Value = 1.8784265
$ws.Range("N20").Value = 3.756853
$ws.Range("O20").Value = 0.08970632314876403
$ws.Range("P20").Value = 0.06164760887349412
$ws.Range("Q20").Value = 0.07343270102233332
$ws.Range("R20").Value = 0.4405962061339999
$ws.Range("S20").Value = 0.0002913606555409777
$ws.Range("T20").Value = 0.0002084547154919064

$ws.Range("G21").Value = 0.03909266666666666
$ws.Range("H21").Value = 0.117278
$ws.Range("I21").Value = 0.00324793888896551
$ws.Range("J21").Value = 0.003381391740913623
$ws.Range("M21").Value = 12.924656
$ws.Range("N21").Value = 38.773968
$ws.Range("O21").Value = 0.617231160081383
$ws.Range("P21").Value = 0.6362565726519981
$ws.Range("Q21").Value = 0.5052592687893332
$ws.Range("R21").Value = 4.547333419104
$ws.Range("S21").Value = 0.00200472908830962
$ws.Range("T21").Value = 0.002151432719867475
